$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 10000
$ws.Range("J21").Value = 10000
$ws.Range("L21").Value = 10000
$ws.Range("N21").Value = -10936
$ws.Range("H23").Value = 10000
$ws.Range("J23").Value = 10000
$ws.Range("L23").Value = 10000
$ws.Range("N23").Value = -10468
$ws.Range("H58").Value = 1732
$ws.Range("I58").Value = 160
$ws.Range("J58").Value = 2125
$ws.Range("K58").Value = 480
$ws.Range("L58").Value = 6375
$ws.Range("M58").Value = -330
$ws.Range("N58").Value = -6675
$ws.Range("H115").Value = 1787.1111
$ws.Range("I115").Value = 760.5
$ws.Range("J115").Value = 10000
$ws.Range("K115").Value = 2281.5
$ws.Range("L115").Value = 30000
$ws.Range("M115").Value = -714.5
$ws.Range("N115").Value = -33134
$ws.Range("H129").Value = 1152.0869
$ws.Range("J129").Value = 1181.7273
$ws.Range("L129").Value = 3545.1819
$ws.Range("N129").Value = -13545.1819
$ws.Range("H133").Value = 42757.777
$ws.Range("J133").Value = 42757.777
$ws.Range("L133").Value = 42757.777
$ws.Range("N133").Value = -52877.777

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22438.803
$ws.Range("I32").Value = 3991.6792
$ws.Range("J32").Value = 348338
$ws.Range("K32").Value = 3991.6792
$ws.Range("L32").Value = 348338
$ws.Range("M32").Value = -3704.6792
$ws.Range("N32").Value = -348912
$ws.Range("H74").Value = 7151.636
$ws.Range("I74").Value = 1312.0769
$ws.Range("J74").Value = 15586.556
$ws.Range("K74").Value = 1312.0769
$ws.Range("L74").Value = 15586.556
$ws.Range("M74").Value = -438.0769
$ws.Range("N74").Value = -17334.556
$ws.Range("H77").Value = 7151.636
$ws.Range("I77").Value = 1312.0769
$ws.Range("J77").Value = 15586.556
$ws.Range("K77").Value = 6560.3845
$ws.Range("L77").Value = 77932.78
$ws.Range("M77").Value = -2192.3845
$ws.Range("N77").Value = -86668.78
$ws.Range("H97").Value = 41679624
$ws.Range("I97").Value = 55572532
$ws.Range("K97").Value = 55572532
$ws.Range("M97").Value = -55572036
$ws.Range("H131").Value = 58466.668
$ws.Range("J131").Value = 58466.668
$ws.Range("L131").Value = 58466.668
$ws.Range("N131").Value = -68546.66800000001
$ws.Range("H133").Value = 63400
$ws.Range("J133").Value = 63400
$ws.Range("L133").Value = 63400
$ws.Range("N133").Value = -68460
$ws.Range("H139").Value = 52966.668
$ws.Range("J139").Value = 52966.668
$ws.Range("L139").Value = 52966.668
$ws.Range("N139").Value = -63246.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 39744
$ws.Range("J130").Value = 39744
$ws.Range("L130").Value = 39744
$ws.Range("N130").Value = -49784
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2098.04
$ws.Range("I58").Value = 1509.5
$ws.Range("J58").Value = 3144.3333
$ws.Range("K58").Value = 1509.5
$ws.Range("L58").Value = 3144.3333
$ws.Range("M58").Value = -1306.5
$ws.Range("N58").Value = -3550.3333
$ws.Range("H132").Value = 3919.72
$ws.Range("I132").Value = 3819.25
$ws.Range("K132").Value = 11457.75
$ws.Range("M132").Value = -8927.75
$ws.Range("H136").Value = 2098.04
$ws.Range("I136").Value = 1509.5
$ws.Range("J136").Value = 3144.3333
$ws.Range("K136").Value = 4528.5
$ws.Range("L136").Value = 9432.999899999999
$ws.Range("M136").Value = -1978.5
$ws.Range("N136").Value = -14532.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 15460.75
$ws.Range("I68").Value = 24486.2
$ws.Range("J68").Value = 418.33334
$ws.Range("K68").Value = 73458.60000000001
$ws.Range("L68").Value = 1255.00002
$ws.Range("M68").Value = -72647.60000000001
$ws.Range("N68").Value = -2877.00002
$ws.Range("H71").Value = 15460.75
$ws.Range("I71").Value = 24486.2
$ws.Range("J71").Value = 418.33334
$ws.Range("K71").Value = 220375.8
$ws.Range("L71").Value = 3765.00006
$ws.Range("M71").Value = -216319.8
$ws.Range("N71").Value = -11877.00006
$ws.Range("H98").Value = 375.375
$ws.Range("I98").Value = 500
$ws.Range("J98").Value = 250.75
$ws.Range("K98").Value = 1500
$ws.Range("L98").Value = 752.25
$ws.Range("M98").Value = -2
$ws.Range("N98").Value = -3748.25
$ws.Range("H103").Value = 3506.25
$ws.Range("I103").Value = 2012.5
$ws.Range("J103").Value = 5000
$ws.Range("K103").Value = 6037.5
$ws.Range("L103").Value = 15000
$ws.Range("M103").Value = -5158.5
$ws.Range("N103").Value = -16758
$ws.Range("H113").Value = 774.95
$ws.Range("I113").Value = 699.5454999999999
$ws.Range("J113").Value = 867.1111
$ws.Range("K113").Value = 2098.6365
$ws.Range("L113").Value = 2601.3333
$ws.Range("M113").Value = 71.36350000000039
$ws.Range("N113").Value = -6941.3333
$ws.Range("H131").Value = 5849303.5
$ws.Range("I131").Value = 392
$ws.Range("J131").Value = 6411699
$ws.Range("K131").Value = 1176
$ws.Range("L131").Value = 19235097
$ws.Range("M131").Value = 3864
$ws.Range("N131").Value = -19245177

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 29995
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()
$ws.Range("H126").Value = 3874.1667
$ws.Range("I126").Value = 2545
$ws.Range("K126").Value = 7635
$ws.Range("M126").Value = -5165
$ws.Range("H132").Value = 1843.5278
$ws.Range("I132").Value = 1754.1666
$ws.Range("K132").Value = 5262.4998
$ws.Range("M132").Value = -2732.4998
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2409.2666
$ws.Range("I61").Value = 3099.5
$ws.Range("K61").Value = 3099.5
$ws.Range("M61").Value = -2897.5
$ws.Range("H68").Value = 1750.75
$ws.Range("I68").Value = 1934.3334
$ws.Range("J68").Value = 1200
$ws.Range("K68").Value = 1934.3334
$ws.Range("L68").Value = 1200
$ws.Range("M68").Value = -1185.3334
$ws.Range("N68").Value = -2698
$ws.Range("H70").Value = 6000
$ws.Range("I70").Value = 6000
$ws.Range("K70").Value = 6000
$ws.Range("M70").Value = -5730
$ws.Range("H71").Value = 1750.75
$ws.Range("I71").Value = 1934.3334
$ws.Range("J71").Value = 1200
$ws.Range("K71").Value = 9671.666999999999
$ws.Range("L71").Value = 6000
$ws.Range("M71").Value = -5927.666999999999
$ws.Range("N71").Value = -13488
$ws.Range("H73").Value = 6000
$ws.Range("I73").Value = 6000
$ws.Range("K73").Value = 6000
$ws.Range("M73").Value = -5064
$ws.Range("H82").Value = 788.75
$ws.Range("I82").Value = 715.7143
$ws.Range("J82").Value = 1300
$ws.Range("K82").Value = 715.7143
$ws.Range("L82").Value = 1300
$ws.Range("M82").Value = -354.7143
$ws.Range("N82").Value = -2022
$ws.Range("H85").Value = 788.75
$ws.Range("I85").Value = 715.7143
$ws.Range("J85").Value = 1300
$ws.Range("K85").Value = 715.7143
$ws.Range("L85").Value = 1300
$ws.Range("M85").Value = 532.2857
$ws.Range("N85").Value = -3796
$ws.Range("H100").Value = 2722.4583
$ws.Range("I100").Value = 1933.9
$ws.Range("J100").Value = 3285.7144
$ws.Range("K100").Value = 1933.9
$ws.Range("L100").Value = 3285.7144
$ws.Range("M100").Value = -1392.9
$ws.Range("N100").Value = -4367.7144
$ws.Range("H113").Value = 2409.2666
$ws.Range("I113").Value = 3099.5
$ws.Range("K113").Value = 3099.5
$ws.Range("M113").Value = -929.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 25001.111
$ws.Range("J123").Value = 25001.111
$ws.Range("L123").Value = 25001.111
$ws.Range("N123").Value = -34801.111
